$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# copy C4's style (fontId2,fillId2,borderId2) fully and then just remove bold? No we need fontId0.
# Let's instead copy fill/border pattern from C4 (border2) but need fillId2 not theme 7. Let's check C4 actual def again: fontId2 fillId2 borderId2. We want fontId0 fillId2 borderId2. So copy C4's Borders + Interior but not font.

$src = $ws.Range("C4")
$dst = $ws.Range("H5")
$dst.Borders.LineStyle = $src.Borders.LineStyle
$dst.Interior.Pattern = $src.Interior.Pattern
$dst.Interior.ThemeColor = $src.Interior.ThemeColor
$dst.HorizontalAlignment = $src.HorizontalAlignment
$dst.VerticalAlignment = $src.VerticalAlignment
"done"
